$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 675.8
$ws.Range("I111").Value = 359.66666
$ws.Range("J111").Value = 1150
$ws.Range("K111").Value = 1078.99998
$ws.Range("L111").Value = 3450
$ws.Range("M111").Value = 1988.00002
$ws.Range("N111").Value = -9584

$ws.Range("H129").Value = 597.0476
$ws.Range("I129").Value = 305.6154
$ws.Range("J129").Value = 1070.625
$ws.Range("K129").Value = 916.8462000000001
$ws.Range("L129").Value = 3211.875
$ws.Range("M129").Value = 4083.1538
$ws.Range("N129").Value = -13211.875

$ws.Range("H137").Value = 2614.3076
$ws.Range("I137").Value = 1927.8833
$ws.Range("J137").Value = 4902.3887
$ws.Range("K137").Value = 5783.6499
$ws.Range("L137").Value = 14707.1661
$ws.Range("M137").Value = -3233.6499
$ws.Range("N137").Value = -19807.1661

$ws.Range("H138").Value = 4778.6353
$ws.Range("I138").Value = 1433.0741
$ws.Range("J138").Value = 9080.071
$ws.Range("K138").Value = 4299.2223
$ws.Range("L138").Value = 27240.213
$ws.Range("M138").Value = 840.7776999999996
$ws.Range("N138").Value = -37520.213

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14246.561
$ws.Range("I32").Value = 3189.5
$ws.Range("J32").Value = 62377.293
$ws.Range("K32").Value = 3189.5
$ws.Range("L32").Value = 62377.293
$ws.Range("M32").Value = -2902.5
$ws.Range("N32").Value = -62951.293

$ws.Range("H61").Value = 1327.4524
$ws.Range("I61").Value = 1072.697
$ws.Range("J61").Value = 2261.5557
$ws.Range("K61").Value = 1072.697
$ws.Range("L61").Value = 2261.5557
$ws.Range("M61").Value = -860.6969999999999
$ws.Range("N61").Value = -2685.5557

$ws.Range("H74").Value = 1475.6923
$ws.Range("I74").Value = 1040
$ws.Range("J74").Value = 1748
$ws.Range("K74").Value = 1040
$ws.Range("L74").Value = 1748
$ws.Range("M74").Value = -166
$ws.Range("N74").Value = -3496

$ws.Range("H77").Value = 1475.6923
$ws.Range("I77").Value = 1040
$ws.Range("J77").Value = 1748
$ws.Range("K77").Value = 5200
$ws.Range("L77").Value = 8740
$ws.Range("M77").Value = -832
$ws.Range("N77").Value = -17476

$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680

$ws.Range("H110").Value = 1033.5714
$ws.Range("I110").Value = 884.4
$ws.Range("J110").Value = 1406.5
$ws.Range("K110").Value = 884.4
$ws.Range("L110").Value = 1406.5
$ws.Range("M110").Value = 1160.6
$ws.Range("N110").Value = -5496.5

$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H132").Value = 2045.1613
$ws.Range("I132").Value = 1373.3585
$ws.Range("J132").Value = 6001.3335
$ws.Range("K132").Value = 4120.0755
$ws.Range("L132").Value = 18004.0005
$ws.Range("M132").Value = -1590.0755
$ws.Range("N132").Value = -23064.0005

$ws.Range("H136").Value = 1327.4524
$ws.Range("I136").Value = 1072.697
$ws.Range("J136").Value = 2261.5557
$ws.Range("K136").Value = 3218.090999999999
$ws.Range("L136").Value = 6784.6671
$ws.Range("M136").Value = -668.0909999999994
$ws.Range("N136").Value = -11884.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1012.0333
$ws.Range("I134").Value = 982.88
$ws.Range("J134").Value = 1157.8
$ws.Range("K134").Value = 2948.64
$ws.Range("L134").Value = 3473.4
$ws.Range("M134").Value = -413.6399999999999
$ws.Range("N134").Value = -8543.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 773.3333
$ws.Range("I16").Value = 683.3333
$ws.Range("J16").Value = 863.3333
$ws.Range("K16").Value = 683.3333
$ws.Range("L16").Value = 863.3333
$ws.Range("M16").Value = -396.3333
$ws.Range("N16").Value = -1437.3333

$ws.Range("H110").Value = 28621.143
$ws.Range("J110").Value = 28621.143
$ws.Range("L110").Value = 28621.143
$ws.Range("N110").Value = -36801.143

$ws.Range("H111").Value = 29498.334
$ws.Range("J111").Value = 29498.334
$ws.Range("L111").Value = 29498.334
$ws.Range("N111").Value = -37678.334

$ws.Range("H113").Value = 773.3333
$ws.Range("I113").Value = 683.3333
$ws.Range("J113").Value = 863.3333
$ws.Range("K113").Value = 683.3333
$ws.Range("L113").Value = 863.3333
$ws.Range("M113").Value = 1486.6667
$ws.Range("N113").Value = -5203.3333

$ws.Range("H114").Value = 28475
$ws.Range("J114").Value = 28475
$ws.Range("L114").Value = 28475
$ws.Range("N114").Value = -37153

$ws.Range("H134").Value = 1156.1774
$ws.Range("I134").Value = 1249.3265
$ws.Range("J134").Value = 805.0769
$ws.Range("K134").Value = 3747.979499999999
$ws.Range("L134").Value = 2415.2307
$ws.Range("M134").Value = -1212.979499999999
$ws.Range("N134").Value = -7485.2307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1582.4062
$ws.Range("I5").Value = 443.3
$ws.Range("J5").Value = 2100.182
$ws.Range("K5").Value = 1329.9
$ws.Range("L5").Value = 6300.545999999999
$ws.Range("M5").Value = -1217.9
$ws.Range("N5").Value = -6524.545999999999

$ws.Range("H113").Value = 476.94736
$ws.Range("I113").Value = 454.22223
$ws.Range("J113").Value = 497.4
$ws.Range("K113").Value = 1362.66669
$ws.Range("L113").Value = 1492.2
$ws.Range("M113").Value = 807.33331
$ws.Range("N113").Value = -5832.2

$ws.Range("H135").Value = 1582.4062
$ws.Range("I135").Value = 443.3
$ws.Range("J135").Value = 2100.182
$ws.Range("K135").Value = 3989.7
$ws.Range("L135").Value = 18901.638
$ws.Range("M135").Value = -1454.7
$ws.Range("N135").Value = -23971.638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4193.7144
$ws.Range("I70").Value = 4154.4614
$ws.Range("J70").Value = 4257.5
$ws.Range("K70").Value = 4154.4614
$ws.Range("L70").Value = 4257.5
$ws.Range("M70").Value = -3884.4614
$ws.Range("N70").Value = -4797.5

$ws.Range("H73").Value = 4193.7144
$ws.Range("I73").Value = 4154.4614
$ws.Range("J73").Value = 4257.5
$ws.Range("K73").Value = 4154.4614
$ws.Range("L73").Value = 4257.5
$ws.Range("M73").Value = -3218.4614
$ws.Range("N73").Value = -6129.5

$ws.Range("H107").Value = 721
$ws.Range("I107").Value = 585.4
$ws.Range("J107").Value = 825.3077
$ws.Range("K107").Value = 585.4
$ws.Range("L107").Value = 825.3077
$ws.Range("M107").Value = 1334.6
$ws.Range("N107").Value = -4665.3077

$ws.Range("H110").Value = 33000
$ws.Range("J110").Value = 33000
$ws.Range("L110").Value = 33000
$ws.Range("N110").Value = -41180

$ws.Range("H111").Value = 25000
$ws.Range("J111").Value = 25000
$ws.Range("L111").Value = 25000
$ws.Range("N111").Value = -31134

$ws.Range("H113").Value = 4250
$ws.Range("I113").Value = 4250
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4250
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2080
$ws.Range("N113").ClearContents()

$ws.Range("H114").Value = 30722
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()

$ws.Range("H122").Value = 2057.0527
$ws.Range("I122").Value = 1685.8667
$ws.Range("J122").Value = 3449
$ws.Range("K122").Value = 5057.6001
$ws.Range("L122").Value = 10347
$ws.Range("M122").Value = -2607.6001
$ws.Range("N122").Value = -15247

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1913.6364
$ws.Range("I61").Value = 1689.7
$ws.Range("J61").Value = 2100.25
$ws.Range("K61").Value = 1689.7
$ws.Range("L61").Value = 2100.25
$ws.Range("M61").Value = -1487.7
$ws.Range("N61").Value = -2504.25

$ws.Range("H110").Value = 29027.5
$ws.Range("J110").Value = 29027.5
$ws.Range("L110").Value = 29027.5
$ws.Range("N110").Value = -37207.5

$ws.Range("H111").Value = 9999.5
$ws.Range("J111").Value = 9999.5
$ws.Range("L111").Value = 9999.5
$ws.Range("N111").Value = -18179.5

$ws.Range("H113").Value = 1913.6364
$ws.Range("I113").Value = 1689.7
$ws.Range("J113").Value = 2100.25
$ws.Range("K113").Value = 1689.7
$ws.Range("L113").Value = 2100.25
$ws.Range("M113").Value = 480.3
$ws.Range("N113").Value = -6440.25

$ws.Range("H114").Value = 30398
$ws.Range("J114").Value = 30398
$ws.Range("L114").Value = 30398
$ws.Range("N114").Value = -39076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 68256.336
$ws.Range("J110").Value = 68256.336
$ws.Range("L110").Value = 68256.336
$ws.Range("N110").Value = -76436.336

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H113").Value = 80357390
$ws.Range("I113").Value = 90909360
$ws.Range("J113").Value = 41666830
$ws.Range("K113").Value = 272728080
$ws.Range("L113").Value = 125000490
$ws.Range("M113").Value = -272725910
$ws.Range("N113").Value = -125004830

$ws.Range("H114").Value = 28000
$ws.Range("J114").Value = 28000
$ws.Range("L114").Value = 28000
$ws.Range("N114").Value = -36678

$ws.Range("H122").Value = 590329.9
$ws.Range("I122").Value = 770992.5600000001
$ws.Range("J122").Value = 3176.25
$ws.Range("K122").Value = 2312977.68
$ws.Range("L122").Value = 9528.75
$ws.Range("M122").Value = -2310527.68
$ws.Range("N122").Value = -14428.75

$ws.Range("H132").Value = 564.1070999999999
$ws.Range("I132").Value = 305.8421
$ws.Range("J132").Value = 1109.3334
$ws.Range("K132").Value = 917.5263
$ws.Range("L132").Value = 3328.0002
$ws.Range("M132").Value = 1612.4737
$ws.Range("N132").Value = -8388.0002

$ws.Range("H136").Value = 364.3846
$ws.Range("I136").Value = 283
$ws.Range("J136").Value = 723.8333
$ws.Range("K136").Value = 849
$ws.Range("L136").Value = 2171.4999
$ws.Range("M136").Value = 1701
$ws.Range("N136").Value = -7271.4999
